# Add a new "footsteps" source entry after the existing
# "334266__projectsu012__short-explosion-1" entry, mirroring the layout
# used by every other entry in the sources list:
#   <name> - <hyperlink to the freesound page (display text == URL)> <space>

$d = $word.ActiveDocument

# The "334266__projectsu012__short-explosion-1" paragraph is the last
# populated entry in the sources list (paragraph 9); the two paragraphs
# that follow it are blank trailer paragraphs. Insert the new entry right
# after paragraph 9, before those trailers.
$anchorRange = $d.Paragraphs(9).Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

$newIndex = 10
$newRange = $d.Paragraphs($newIndex).Range
$newRange.Collapse(1)

$title = "479445__yatoimtop__footsteps"
$sep = " - "
$url = "https://freesound.org/people/yatoimtop/sounds/479445/"

# Type the plain-text title, separator and URL first ...
$newRange.InsertAfter($title + $sep + $url)

# ... then convert the trailing URL text into a real hyperlink (display
# text equal to the address), just like the rest of the document.
$fullRange = $d.Paragraphs($newIndex).Range
$urlStart = $fullRange.Start + ($title + $sep).Length
$urlEnd = $fullRange.End
$urlRange = $d.Range($urlStart, $urlEnd)

[void]$d.Hyperlinks.Add($urlRange, $url, "", "", $url)

# Finally append the trailing space that terminates every entry.
$urlRange.InsertAfter(" ")
